$d = $word.ActiveDocument

# Hybrid bold + color (2C3E50) highlighting for quantitative impact metrics.
# Color 2C3E50 (R=0x2C,G=0x3E,B=0x50) expressed as the Word/VBA BGR-packed
# long that RGB() produces: R + G*256 + B*65536 = 44 + 62*256 + 80*65536
$metricColor = 5258796

function Highlight {
    param($ParaIndex, $SearchText, $StartPos)
    $para = $d.Paragraphs.Item($ParaIndex)
    $rng = $para.Range
    if ($StartPos -ne -1) {
        $rng.Start = $StartPos
    }
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $metricColor
    }
    return $rng.End
}

# --- Paragraph 10: "...demographic classification accuracy from 23% to 64%"
$pos = Highlight 10 "23%" -1
$pos = Highlight 10 "64%" $pos

# --- Paragraph 12: "...margin of error from ±4.2% to ±2.1%, increasing voter turnout
#     prediction accuracy from 71% to 87%, and ensuring..."
$pos = Highlight 12 "±4.2%" -1
$pos = Highlight 12 "±2.1%" $pos
$pos = Highlight 12 "71%" $pos
$pos = Highlight 12 "87%" $pos

# --- Paragraph 13: "...reduced mapping costs by 73.5%, saving campaigns and
#     organizations $4.7M and enabling..."
$pos = Highlight 13 "73.5%" -1
$pos = Highlight 13 "$4.7M" $pos

# --- Paragraph 14: "...valued over $2 trillion"
$pos = Highlight 14 "$2" -1

# --- Paragraph 39: "...reducing processing time by 57%"
$pos = Highlight 39 "57%" -1

# --- Paragraph 55: "Algorithmic innovation: ... reducing mapping costs 73.5%"
$pos = Highlight 55 "73.5%" -1

# --- Paragraph 56: "$4.7M savings enabled nonprofit access"
$pos = Highlight 56 "$4.7M" -1

# --- Paragraph 57: "Platform impact: ... serving 12,847 analysts across 89 organizations"
$pos = Highlight 57 "12,847" -1
